$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1300.5454
$ws.Range("I129").Value = 291.8
$ws.Range("J129").Value = 2141.1667
$ws.Range("K129").Value = 875.4000000000001
$ws.Range("L129").Value = 6423.500100000001
$ws.Range("M129").Value = 4124.6
$ws.Range("N129").Value = -16423.5001
$ws.Range("H137").Value = 24391910
$ws.Range("I137").Value = 30303944
$ws.Range("J137").Value = 4775
$ws.Range("K137").Value = 90911832
$ws.Range("L137").Value = 14325
$ws.Range("M137").Value = -90909282
$ws.Range("N137").Value = -19425
$ws.Range("H141").Value = 1999.762
$ws.Range("I141").Value = 1999.762
$ws.Range("K141").Value = 5999.286
$ws.Range("M141").Value = -819.2860000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21543.56
$ws.Range("I32").Value = 6191.627
$ws.Range("J32").Value = 150116
$ws.Range("K32").Value = 6191.627
$ws.Range("L32").Value = 150116
$ws.Range("M32").Value = -5904.627
$ws.Range("N32").Value = -150690
$ws.Range("H45").Value = 1287.6316
$ws.Range("I45").Value = 1292.0588
$ws.Range("J45").Value = 1250
$ws.Range("K45").Value = 1292.0588
$ws.Range("L45").Value = 1250
$ws.Range("M45").Value = -915.0588
$ws.Range("N45").Value = -2004
$ws.Range("H74").Value = 3656.0815
$ws.Range("I74").Value = 814.6667
$ws.Range("J74").Value = 9516.5
$ws.Range("K74").Value = 814.6667
$ws.Range("L74").Value = 9516.5
$ws.Range("M74").Value = 59.33330000000001
$ws.Range("N74").Value = -11264.5
$ws.Range("H77").Value = 3656.0815
$ws.Range("I77").Value = 814.6667
$ws.Range("J77").Value = 9516.5
$ws.Range("K77").Value = 4073.3335
$ws.Range("L77").Value = 47582.5
$ws.Range("M77").Value = 294.6665000000003
$ws.Range("N77").Value = -56318.5
$ws.Range("H88").Value = 2810
$ws.Range("I88").Value = 5700
$ws.Range("J88").Value = 2087.5
$ws.Range("K88").Value = 5700
$ws.Range("L88").Value = 2087.5
$ws.Range("M88").Value = -5294
$ws.Range("N88").Value = -2899.5
$ws.Range("H91").Value = 2810
$ws.Range("I91").Value = 5700
$ws.Range("J91").Value = 2087.5
$ws.Range("K91").Value = 5700
$ws.Range("L91").Value = 2087.5
$ws.Range("M91").Value = -4296
$ws.Range("N91").Value = -4895.5
$ws.Range("H97").Value = 5504.85
$ws.Range("I97").Value = 5741.9473
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 5741.9473
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -5245.9473
$ws.Range("N97").Value = -1992
$ws.Range("H122").Value = 3655.5
$ws.Range("I122").Value = 3719.9333
$ws.Range("K122").Value = 11159.7999
$ws.Range("M122").Value = -8709.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -21232
$ws.Range("H94").Value = 2712.2964
$ws.Range("I94").Value = 2314.4348
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 2314.4348
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -1863.4348
$ws.Range("N94").Value = -5902
$ws.Range("H99").Value = 1947.6666
$ws.Range("I99").Value = 1947.6666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1947.6666
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -449.6666
$ws.Range("N99").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H86").Value = 125002620
$ws.Range("I86").Value = 250000750
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 250000750
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -249999627
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 125002620
$ws.Range("I89").Value = 250000750
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 1250003750
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -1249998134
$ws.Range("N89").Value = -33732
$ws.Range("H105").Value = 928
$ws.Range("I105").Value = 812.82355
$ws.Range("J105").Value = 1078.6154
$ws.Range("K105").Value = 812.82355
$ws.Range("L105").Value = 1078.6154
$ws.Range("M105").Value = 934.17645
$ws.Range("N105").Value = -4572.6154

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 144
$ws.Range("I36").Value = 144
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 432
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -263
$ws.Range("N36").Value = ""
$ws.Range("H123").Value = 2791
$ws.Range("I123").Value = 1515
$ws.Range("J123").Value = 3216.3333
$ws.Range("K123").Value = 4545
$ws.Range("L123").Value = 9648.999899999999
$ws.Range("M123").Value = -2095
$ws.Range("N123").Value = -14548.9999
$ws.Range("H124").Value = 250000860
$ws.Range("I124").Value = 915
$ws.Range("K124").Value = 2745
$ws.Range("M124").Value = 2165
$ws.Range("H131").Value = 7577283
$ws.Range("J131").Value = 8548625
$ws.Range("L131").Value = 25645875
$ws.Range("N131").Value = -25655955
$ws.Range("H134").Value = 7654
$ws.Range("I134").Value = 3928.7856
$ws.Range("J134").Value = 10398.895
$ws.Range("K134").Value = 11786.3568
$ws.Range("L134").Value = 31196.685
$ws.Range("M134").Value = -6716.356800000001
$ws.Range("N134").Value = -41336.685
$ws.Range("H137").Value = 6317502
$ws.Range("I137").Value = 10004896
$ws.Range("J137").Value = 171844.33
$ws.Range("K137").Value = 30014688
$ws.Range("L137").Value = 515532.99
$ws.Range("M137").Value = -30009588
$ws.Range("N137").Value = -525732.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1254.0952
$ws.Range("I97").Value = 1254.5264
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 1254.5264
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -758.5264
$ws.Range("N97").Value = -2242
$ws.Range("H132").Value = 3056.1667
$ws.Range("I132").Value = 2606.8696
$ws.Range("J132").Value = 3851.077
$ws.Range("K132").Value = 7820.6088
$ws.Range("L132").Value = 11553.231
$ws.Range("M132").Value = -5290.6088
$ws.Range("N132").Value = -16613.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 677.2857
$ws.Range("I46").Value = 697
$ws.Range("J46").Value = 662.5
$ws.Range("K46").Value = 697
$ws.Range("L46").Value = 662.5
$ws.Range("M46").Value = -509
$ws.Range("N46").Value = -1038.5
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = ""
$ws.Range("H132").Value = 3619.375
$ws.Range("I132").Value = 2625.6667
$ws.Range("J132").Value = 4897
$ws.Range("K132").Value = 7877.000100000001
$ws.Range("L132").Value = 14691
$ws.Range("M132").Value = -5347.000100000001
$ws.Range("N132").Value = -19751
$ws.Range("H136").Value = 7047.7095
$ws.Range("I136").Value = 5741.6113
$ws.Range("J136").Value = 8856.154
$ws.Range("K136").Value = 17224.8339
$ws.Range("L136").Value = 26568.462
$ws.Range("M136").Value = -14674.8339
$ws.Range("N136").Value = -31668.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 500002200
$ws.Range("I96").Value = 1000000000
$ws.Range("J96").Value = 4444
$ws.Range("K96").Value = 1000000000
$ws.Range("L96").Value = 4444
$ws.Range("M96").Value = -999998627
$ws.Range("N96").Value = -7190
$ws.Range("H107").Value = 11112612
$ws.Range("I107").Value = 11112612
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 33337836
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -33335916
$ws.Range("N107").Value = ""
$ws.Range("H136").Value = 3050.12
$ws.Range("I136").Value = 815.3913
$ws.Range("J136").Value = 28749.5
$ws.Range("K136").Value = 2446.1739
$ws.Range("L136").Value = 86248.5
$ws.Range("M136").Value = 103.8261000000002
$ws.Range("N136").Value = -91348.5
